# Insert a new weekly record at row 589 (pushes existing rows 589-634 down to 590-635)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(589).Insert()

$ws.Cells.Item(589, 1).Value = 10
$ws.Cells.Item(589, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(589, 3).Value = "La Araucanía"
$ws.Cells.Item(589, 4).Value = 45021
$ws.Cells.Item(589, 5).Value = 9
$ws.Cells.Item(589, 6).Value = 100112023
$ws.Cells.Item(589, 7).Value = "Brócoli"
$ws.Cells.Item(589, 8).Value = "Sin especificar"
$ws.Cells.Item(589, 9).Value = "Primera"
$ws.Cells.Item(589, 10).Value = 2100
$ws.Cells.Item(589, 11).Value = 1200
$ws.Cells.Item(589, 12).Value = 1200
$ws.Cells.Item(589, 13).Value = 1200
$ws.Cells.Item(589, 14).Value = "`$/unidad"
$ws.Cells.Item(589, 15).Value = "Región del Maule"
$ws.Cells.Item(589, 16).Value = 1200
$ws.Cells.Item(589, 17).Value = 1
$ws.Cells.Item(589, 18).Value = "Hortaliza"

$ws.Cells.Item(589, 4).NumberFormat = $ws.Cells.Item(590, 4).NumberFormat
